$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Avals = @("ECs","ECs","ECs","ECs","ECs","Inflammatory-Mac","Inflammatory-Mac","Inflammatory-Mac","Inflammatory-Mac","Inflammatory-Mac","Resolving-Mac","Resolving-Mac","Resolving-Mac","Resolving-Mac","Resolving-Mac")
$Bvals = @("Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk","Ifnk")
$Cvals = @("Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2","Ifnar2")
$Dvals = @("ECs","FAPs","Inflammatory-Mac","MuSCs","Resolving-Mac","ECs","FAPs","Inflammatory-Mac","MuSCs","Resolving-Mac","ECs","FAPs","Inflammatory-Mac","MuSCs","Resolving-Mac")
$Evals = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$Fvals = @(0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333,0.3333333333333333)
$Gvals = @(0.2027413333333333,0.2027413333333333,0.2027413333333333,0.2027413333333333,0.2027413333333333,0.7050673333333334,0.7050673333333334,0.7050673333333334,0.7050673333333334,0.7050673333333334,0.152942,0.152942,0.152942,0.152942,0.152942)
$Hvals = @(0.608224,0.608224,0.608224,0.608224,0.608224,2.115202,2.115202,2.115202,2.115202,2.115202,0.458826,0.458826,0.458826,0.458826,0.458826)
$Ivals = @(0.1911300550679205,0.1911300550679205,0.1911300550679205,0.1911300550679205,0.1911300550679205,0.6646871460839683,0.6646871460839683,0.6646871460839683,0.6646871460839683,0.6646871460839683,0.1441827988481113,0.1441827988481113,0.1441827988481113,0.1441827988481113,0.1441827988481113)
$Jvals = @(0.1911300550679204,0.1911300550679204,0.1911300550679204,0.1911300550679204,0.1911300550679204,0.6646871460839682,0.6646871460839682,0.6646871460839682,0.6646871460839682,0.6646871460839682,0.1441827988481113,0.1441827988481113,0.1441827988481113,0.1441827988481113,0.1441827988481113)
$Kvals = @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3)
$Lvals = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$Mvals = @(27.85292233333333,26.66082666666667,105.665011,4.883238666666666,71.84303666666666,27.85292233333333,26.66082666666667,105.665011,4.883238666666666,71.84303666666666,27.85292233333333,26.66082666666667,105.665011,4.883238666666666,71.84303666666666)
$Nvals = @(83.55876699999999,79.98248000000001,316.995033,14.649716,215.52911,83.55876699999999,79.98248000000001,316.995033,14.649716,215.52911,83.55876699999999,79.98248000000001,316.995033,14.649716,215.52911)
$Ovals = @(0.1175699887262562,0.1125380329259528,0.4460226472237104,0.02061264193813266,0.3032566891859479,0.1175699887262562,0.1125380329259528,0.4460226472237104,0.02061264193813266,0.3032566891859479,0.1175699887262562,0.1125380329259528,0.4460226472237104,0.02061264193813266,0.3032566891859479)
$Pvals = @(0.1175699887262562,0.1125380329259528,0.4460226472237104,0.02061264193813266,0.3032566891859479,0.1175699887262562,0.1125380329259528,0.4460226472237104,0.02061264193813266,0.3032566891859479,0.1175699887262562,0.1125380329259528,0.4460226472237104,0.02061264193813266,0.3032566891859479)
$Qvals = @(5.646938611089777,5.405251546168889,21.42266521682133,0.9900343182648887,14.56555304451556,19.63818567510378,18.79767796232889,74.50094753240734,3.443012064736889,50.65417828113556,4.259881647504666,4.077560152053334,16.160618112362,0.7468522881573333,10.98781771387333)
$Rvals = @(50.82244749980799,48.64726391552001,192.803986951392,8.910308864384,131.08997740064,176.743671075934,169.17910166096,670.508527791666,30.987108582632,455.88760453022,38.338934827542,36.69804136848001,145.445563011258,6.721670593416,98.89035942486001)
$Svals = @(0.02247115841958414,0.0215094004303728,0.08524833312540744,0.003939695388730622,0.05796146770382546,0.07814726027157956,0.0748025839314552,0.2964655204719447,0.01370095814310811,0.2015708232658807,0.01695157003509251,0.01622604856412477,0.06430879362635836,0.002971988406293925,0.04372439821624174)
$Tvals = @(0.02247115841958413,0.0215094004303728,0.08524833312540743,0.003939695388730622,0.05796146770382545,0.07814726027157955,0.07480258393145518,0.2964655204719446,0.01370095814310811,0.2015708232658807,0.01695157003509251,0.01622604856412478,0.06430879362635836,0.002971988406293926,0.04372439821624174)

for ($i = 0; $i -lt 15; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $Avals[$i]
    $ws.Cells.Item($r, 2).Value = $Bvals[$i]
    $ws.Cells.Item($r, 3).Value = $Cvals[$i]
    $ws.Cells.Item($r, 4).Value = $Dvals[$i]
    $ws.Cells.Item($r, 5).Value = $Evals[$i]
    $ws.Cells.Item($r, 6).Value = $Fvals[$i]
    $ws.Cells.Item($r, 7).Value = $Gvals[$i]
    $ws.Cells.Item($r, 8).Value = $Hvals[$i]
    $ws.Cells.Item($r, 9).Value = $Ivals[$i]
    $ws.Cells.Item($r, 10).Value = $Jvals[$i]
    $ws.Cells.Item($r, 11).Value = $Kvals[$i]
    $ws.Cells.Item($r, 12).Value = $Lvals[$i]
    $ws.Cells.Item($r, 13).Value = $Mvals[$i]
    $ws.Cells.Item($r, 14).Value = $Nvals[$i]
    $ws.Cells.Item($r, 15).Value = $Ovals[$i]
    $ws.Cells.Item($r, 16).Value = $Pvals[$i]
    $ws.Cells.Item($r, 17).Value = $Qvals[$i]
    $ws.Cells.Item($r, 18).Value = $Rvals[$i]
    $ws.Cells.Item($r, 19).Value = $Svals[$i]
    $ws.Cells.Item($r, 20).Value = $Tvals[$i]
}
